# iOS ratings workbook update:
#  - Insert a new "Detail Date" column (B) with a precise timestamp,
#    shifting the existing Date/App Name/Rating/Reviews/Rank columns one
#    to the right (B->C, C->D, D->E, E->F, F->G).
#  - Refresh the iOS Total Reviews (now column F) and iOS Rank (now
#    column G) figures with the latest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "Detail Date" column before the existing Date column ---
$ws.Columns.Item(2).Insert()

# Header cell: give it the same (bold / bordered) header style as the
# other header cells by copying formats from the neighboring header cell.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B1").Value = "Detail Date"

# Data cells B2:B38 inherited column A's style (bold/border) on insert;
# clear that so they match the plain formatting used by the rest of the
# data rows, then fill in the detail timestamp.
$detailRange = $ws.Range("B2:B38")
$detailRange.ClearFormats()
$detailRange.Value = "2023-06-27 19:52:58"

# --- Refresh iOS Total Reviews (F) and iOS Rank (G) values ---
$reviews = @{
    2  = 61782
    3  = 4429
    4  = 123
    5  = 399238
    6  = 4754964
    7  = 43165
    8  = 322313
    9  = 2412815
    10 = 42797
    11 = 24983
    12 = 2011063
    13 = 1047295
    14 = 258123
    15 = 116079
    16 = 77175
    17 = 175
    18 = 411836
    19 = 3991292
    20 = 4150
    21 = 6605
    22 = 3129
    23 = 45223
    24 = 877477
    25 = 298
    26 = 1623
    27 = 14
    28 = 1
    29 = 307
    30 = 20
    31 = 24
    32 = 1375
    33 = 28390
    34 = 8
    35 = 13445
    36 = 38
    37 = 3300
    38 = 22
}

$rank = @{
    5  = 104
    6  = 9
    7  = 91
    9  = 4
    10 = 89
    11 = 64
    12 = 10
    14 = 97
    15 = 171
    16 = 141
    18 = 50
    19 = 11
    22 = 68
    23 = 178
    24 = 5
}

for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 6).Value = $reviews[$r]
    if ($rank.ContainsKey($r)) {
        $ws.Cells.Item($r, 7).Value = $rank[$r]
    }
}
